$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.048.74"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.302.36"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.25%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.94"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "2.660.72"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "2.282.19"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "42.927.76"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.68%  "
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0687"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.111"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").Value = "2.012.56"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0288"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "2.527.51"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.83%  "
